$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.169.05"
$ws.Range("D3").Value = "3.106.53"
$ws.Range("E3").Value = "  -0.20%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "578.77"
$ws.Range("E5").Value = "  -0.12%  "
$ws.Range("D6").Value = "173.16"
$ws.Range("E6").Value = "  +0.05%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("E8").Value = "  -0.89%  "
$ws.Range("D9").Value = "6.51"
$ws.Range("E9").Value = "  +1.10%  "
$ws.Range("E10").Value = "  -1.65%  "
$ws.Range("D11").Value = "0.476"
$ws.Range("E11").Value = "  -1.22%  "
$ws.Range("E12").Value = "  -1.07%  "
$ws.Range("D13").Value = "36.71"
$ws.Range("E13").Value = "  -1.46%  "
$ws.Range("E14").Value = "  -1.77%  "
$ws.Range("D15").Value = "3.624.04"
$ws.Range("E15").Value = "  -0.15%  "
$ws.Range("D16").Value = "67.125.84"
$ws.Range("E16").Value = "  +0.17%  "
$ws.Range("E17").Value = "  -1.74%  "
$ws.Range("D18").Value = "3.111.75"
$ws.Range("E18").Value = "  -0.04%  "
$ws.Range("D19").Value = "16.51"
$ws.Range("E19").Value = "  +1.68%  "
$ws.Range("D20").Value = "490.66"
$ws.Range("E20").Value = "  +0.71%  "
$ws.Range("D21").Value = "0.701"
$ws.Range("E21").Value = "  -2.46%  "
$ws.Range("D22").Value = "7.83"
$ws.Range("E22").Value = "  +2.93%  "
$ws.Range("D23").Value = "83.90"
$ws.Range("E23").Value = "  -0.86%  "
$ws.Range("D24").Value = "13.05"
$ws.Range("E24").Value = "  -2.18%  "
$ws.Range("D25").Value = "2.28"
$ws.Range("E25").Value = "  -3.65%  "
$ws.Range("D26").Value = "10.57"
$ws.Range("E26").Value = "  +4.81%  "
$ws.Range("E27").Value = "  -0.08%  "
$ws.Range("D28").Value = "7.90"
$ws.Range("E28").Value = "  -2.06%  "
$ws.Range("E29").Value = "  -3.53%  "
$ws.Range("E30").Value = "  -1.15%  "
$ws.Range("D31").Value = "28.25"
$ws.Range("E31").Value = "  -2.66%  "
$ws.Range("E32").Value = "  -1.27%  "
$ws.Range("D33").Value = "0.0₃0928"
$ws.Range("E33").Value = "  -7.37%  "
$ws.Range("D34").Value = "0.999"
$ws.Range("E34").Value = "  -0.08%  "
$ws.Range("D35").Value = "5.82"
$ws.Range("E35").Value = "  -2.20%  "
$ws.Range("E36").Value = "  -1.82%  "
$ws.Range("E37").Value = "  -0.66%  "
$ws.Range("D38").Value = "2.03"
$ws.Range("E38").Value = "  -3.86%  "
$ws.Range("D39").Value = "0.307"
$ws.Range("E39").Value = "  -2.74%  "
$ws.Range("E40").Value = "  +0.46%  "
$ws.Range("D41").Value = "8.45"
$ws.Range("E41").Value = "  -2.71%  "
$ws.Range("D42").Value = "385.21"
$ws.Range("E42").Value = "  -0.07%  "
$ws.Range("D43").Value = "2.797.83"
$ws.Range("E43").Value = "  -1.79%  "
$ws.Range("D44").Value = "2.55"
$ws.Range("E44").Value = "  -8.84%  "
$ws.Range("D45").Value = "0.0350"
$ws.Range("E45").Value = "  -2.44%  "
$ws.Range("D46").Value = "134.98"
$ws.Range("E46").Value = "  -1.35%  "
$ws.Range("D48").Value = "24.96"
$ws.Range("E48").Value = "  -0.78%  "
$ws.Range("D49").Value = "2.19"
$ws.Range("E49").Value = "  -1.83%  "
$ws.Range("E50").Value = "  -1.57%  "
$ws.Range("D51").Value = "6.68"
$ws.Range("E51").Value = "  -2.89%  "
